# Implemented feature no.7 - user permission management.
#
# "product backlog" (sheet 1): user story #7 ("User permissions management")
# and user story #9 ("login and use the application only if I am activated",
# a sub-task of the same feature) are marked as Completed. Story #7's
# open-ended "I want to" note is cleared now that the feature is fully
# specified via the linked product doc.
#
# "sprint backlog" (sheet 2): the now-finished sprint items for stories #7
# and #9 are removed, leaving only the still-open item (#6).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "product backlog"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Row 7 (user story #7 "User permissions management"): the freeform
# "I want to" note is superseded by the status update; mark it Completed.
$ws1.Range("E7").ClearContents()
$ws1.Range("I7").Value = "(3) Completed"

# Row 9 (user story #9, related sub-task): mark it Completed too.
$ws1.Range("I9").Value = "(3) Completed"

# Re-style both rows like the other "Completed"/highlighted rows (2 and 4):
# yellow fill + wrap text across the full A:I span, row height 30.
$rng7 = $ws1.Range("A7:I7")
$rng7.Interior.Color = 65535
$rng7.WrapText = $true
$ws1.Rows.Item(7).RowHeight = 30

$rng9 = $ws1.Range("A9:I9")
$rng9.Interior.Color = 65535
$rng9.WrapText = $true
$ws1.Rows.Item(9).RowHeight = 30

# ---------------------------------------------------------------------
# Sheet 2: "sprint backlog"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Remove the two sprint-backlog rows belonging to the now-completed
# stories (#9 and #7); the remaining row (#6) shifts up to row 2.
$ws2.Rows.Item(2).Delete()
$ws2.Rows.Item(2).Delete()

# Deleting rows shrinks the whole-column data validations (which extended
# to the sheet's last row) by the same amount; restore them to span the
# full column again, in their original order (I, then H).
$iVal = $ws2.Range("I1:I1048576")
$iVal.Validation.Delete()
$iVal.Validation.Add(3, 1, 1, "userstorystatus")

$hVal = $ws2.Range("H1:H1048576")
$hVal.Validation.Delete()
$hVal.Validation.Add(3, 1, 1, "priority")

# Restore the selection to where the user last left it while trimming the
# sheet (cursor now rests past the shrunk data range, on E4).
$ws2.Activate()
$ws2.Range("E4").Select()
